$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "57.382.43"
Set-TextValue "E2" "  -4.64%  "
Set-TextValue "D3" "2.929.24"
Set-TextValue "E3" "  -0.88%  "
Set-TextValue "E4" "  +0.28%  "
Set-TextValue "D5" "549.82"
Set-TextValue "E5" "  -3.80%  "
Set-TextValue "D6" "129.45"
Set-TextValue "E6" "  +5.61%  "
Set-TextValue "E7" "  +0.29%  "
Set-TextValue "D8" "0.514"
Set-TextValue "E8" "  +3.67%  "
Set-TextValue "D9" "2.923.24"
Set-TextValue "E9" "  -1.00%  "
Set-TextValue "D10" "0.128"
Set-TextValue "E10" "  -2.68%  "
Set-TextValue "D11" "4.78"
Set-TextValue "E11" "  -4.91%  "
Set-TextValue "D12" "0.442"
Set-TextValue "E12" "  +1.98%  "
Set-TextValue "D13" "0.0000220"
Set-TextValue "E13" "  +0.00%  "
Set-TextValue "D14" "32.36"
Set-TextValue "E14" "  +0.14%  "
Set-TextValue "D15" "0.120"
Set-TextValue "E15" "  +1.57%  "
Set-TextValue "D16" "3.426.29"
Set-TextValue "E16" "  -0.30%  "
Set-TextValue "D17" "6.71"
Set-TextValue "E17" "  +9.58%  "
Set-TextValue "D18" "2.939.99"
Set-TextValue "E18" "  -0.40%  "
Set-TextValue "D19" "57.613.91"
Set-TextValue "E19" "  -4.21%  "
Set-TextValue "D20" "414.25"
Set-TextValue "E20" "  -3.06%  "
Set-TextValue "D21" "12.99"
Set-TextValue "E21" "  +0.17%  "
Set-TextValue "D22" "0.676"
Set-TextValue "E22" "  +3.23%  "
Set-TextValue "D23" "6.90"
Set-TextValue "E23" "  -0.38%  "
Set-TextValue "D24" "12.87"
Set-TextValue "E24" "  +0.91%  "
Set-TextValue "D25" "78.69"
Set-TextValue "E25" "  +0.73%  "
Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  -0.11%  "
Set-TextValue "E27" "  +0.63%  "
Set-TextValue "D28" "2.47"
Set-TextValue "E28" "  -1.05%  "
Set-TextValue "D29" "7.44"
Set-TextValue "E29" "  +4.90%  "
Set-TextValue "D30" "1.97"
Set-TextValue "E30" "  +5.95%  "
Set-TextValue "D31" "6.12"
Set-TextValue "E31" "  +1.30%  "
Set-TextValue "D32" "0.103"
Set-TextValue "E32" "  +12.47%  "
Set-TextValue "D33" "24.85"
Set-TextValue "E33" "  -0.85%  "
Set-TextValue "D34" "5.58"
Set-TextValue "E34" "  +1.12%  "
Set-TextValue "D35" "2.09"
Set-TextValue "E35" "  -3.61%  "
Set-TextValue "D36" "0.930"
Set-TextValue "E36" "  -0.93%  "
Set-TextValue "D37" "48.33"
Set-TextValue "E37" "  -2.12%  "
Set-TextValue "D38" "0.0₃0672"
Set-TextValue "E38" "  +4.42%  "
Set-TextValue "D39" "8.33"
Set-TextValue "E39" "  +6.64%  "
Set-TextValue "D40" "2.55"
Set-TextValue "E40" "  +7.14%  "
Set-TextValue "D41" "0.0347"
Set-TextValue "E41" "  -1.75%  "
Set-TextValue "D42" "0.107"
Set-TextValue "E42" "  -0.58%  "
Set-TextValue "D43" "376.13"
Set-TextValue "E43" "  +0.72%  "
Set-TextValue "D44" "2.635.81"
Set-TextValue "E44" "  +1.09%  "
Set-TextValue "E45" "  -0.02%  "
Set-TextValue "D46" "0.238"
Set-TextValue "E46" "  +1.89%  "
Set-TextValue "D47" "120.77"
Set-TextValue "E47" "  +1.49%  "
Set-TextValue "E48" "  +3.05%  "
Set-TextValue "D49" "1.97"
Set-TextValue "E49" "  +1.27%  "
Set-TextValue "D50" "23.15"
Set-TextValue "E50" "  +0.27%  "
Set-TextValue "D51" "1.98"
Set-TextValue "E51" "  +1.21%  "
